$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.694.49"
$ws.Range("E2").Value = "  -1.94%  "

$ws.Range("D3").Value = "2.412.09"
$ws.Range("E3").Value = "  +4.95%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'300.39"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").Value = "'97.93"
$ws.Range("E6").Value = "  -2.74%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  -2.02%  "

$ws.Range("D10").Value = "'35.01"
$ws.Range("E10").Value = "  -4.29%  "

$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").Value = "2.765.40"
$ws.Range("E14").Value = "  +4.36%  "

$ws.Range("D15").Value = "2.416.60"
$ws.Range("E15").Value = "  +5.05%  "

$ws.Range("D16").Value = "'0.847"
$ws.Range("E16").Value = "  +4.02%  "

$ws.Range("D17").Value = "'14.28"
$ws.Range("E17").Value = "  +3.10%  "

$ws.Range("D18").Value = "45.679.57"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("D19").Value = "'13.16"
$ws.Range("E19").Value = "  +0.82%  "

$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("D22").Value = "'67.34"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").Value = "'243.65"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("D24").Value = "'2.80"
$ws.Range("E24").Value = "  -3.24%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'1.94"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("D27").Value = "'38.82"
$ws.Range("E27").Value = "  -9.17%  "

$ws.Range("E28").Value = "  -1.97%  "

$ws.Range("D29").Value = "'9.82"
$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("D30").Value = "'3.84"
$ws.Range("E30").Value = "  +16.71%  "

$ws.Range("D31").Value = "'21.35"
$ws.Range("E31").Value = "  +6.73%  "

$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'2.75"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.57"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("D34").Value = "'148.22"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").Value = "'0.0777"
$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("D36").Value = "'1.99"
$ws.Range("E36").Value = "  +12.04%  "

$ws.Range("D37").Value = "'0.113"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("D39").Value = "'15.23"
$ws.Range("E39").Value = "  -5.29%  "

$ws.Range("D40").Value = "'3.91"
$ws.Range("E40").Value = "  -2.17%  "

$ws.Range("D41").Value = "'0.0300"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").Value = "'3.28"
$ws.Range("E42").Value = "  -2.30%  "

$ws.Range("D43").Value = "1.954.48"
$ws.Range("E43").Value = "  +7.38%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "'91.65"
$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("D46").Value = "'1.76"
$ws.Range("E46").Value = "  -10.91%  "

$ws.Range("D47").Value = "'8.71"
$ws.Range("E47").Value = "  +10.28%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'15.52"
$ws.Range("E48").Value = "  +16.56%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'102.38"
$ws.Range("E49").Value = "  +6.54%  "

$ws.Range("D50").Value = "'0.188"
$ws.Range("E50").Value = "  -3.46%  "

$ws.Range("D51").Value = "2.649.78"
$ws.Range("E51").Value = "  +4.91%  "
